$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old helper cells D2/E2 (their only use was feeding D4's formula) ---
$ws.Range("D2:E2").ClearContents()

# --- D4 used to be "=D2*E2"; bake it down to a literal value ---
$ws.Range("D4").Value = 6000000

# --- A5 label "V" becomes "Volume" ---
$ws.Range("A5").Value = "Volume"

# --- New italic note above row 8, in the previously-empty row 7 ---
$ws.Range("F7").Value = "cfast temperature and pressur columns are copied from a CFAST run"
$ws.Range("F7").Font.Italic = $true

# --- New bold column header F9 ---
$ws.Range("F9").Value = "temperature"

# --- C8 label reworded ---
$ws.Range("C8").Value = "calculated pressure"

# --- A9 label "T" becomes "Time" ---
$ws.Range("A9").Value = "Time"

# --- Italic annotation cells trailing rows 9-13 (N:Q), text only in column N ---
$ws.Range("N10").Value = "DP=(gamma-1)*qtotal*Time/Volume"
$ws.Range("N11").Value = "M=M0+mfire*Time"
$ws.Range("N12").Value = "E=E0+qconvec*Time"
$ws.Range("N13").Value = "T=E/(cv*M)-273.3"

# --- New italic annotation block to the right of the table (row 8) ---
$ws.Range("N8").Value = "Formulas (assuming constant fire)"

$ws.Range("N8:Q8").Font.Italic = $true
$ws.Range("N9:Q9").Font.Italic = $true
$ws.Range("N10:Q10").Font.Italic = $true
$ws.Range("N11:Q11").Font.Italic = $true
$ws.Range("N12:Q12").Font.Italic = $true
$ws.Range("N13:Q13").Font.Italic = $true

# --- Update the saved cursor position to match the authored file ---
[void]$ws.Range("D5").Select()
